$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @{
    14 = @(" Dubai (DSC)", " October 25 2020", "Super Kings won by 8 wickets (with 8 balls remaining)", "Royal Challengers Bangalore", "Chennai Super Kings", "AB de Villiers †", "39", "36", "4", "0", "108.33")
    15 = @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", "AB de Villiers †", "15", "12", "1", "1", "125.00")
    16 = @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", "AB de Villiers †", "35", "21", "1", "2", "166.66")
    17 = @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "AB de Villiers †", "56", "43", "5", "0", "130.23")
    18 = @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "AB de Villiers †", "24", "24", "1", "1", "100.00")
    19 = @(" Sharjah", " October 15 2020", "Kings XI won by 8 wickets", "Royal Challengers Bangalore", "Kings XI Punjab", "AB de Villiers †", "2", "5", "0", "0", "40.00")
    20 = @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", "AB de Villiers †", "9", "6", "2", "0", "150.00")
    21 = @(" Dubai (DSC)", " September 28 2020", "Match tied (RCB won the one-over eliminator)", "Royal Challengers Bangalore", "Mumbai Indians", "AB de Villiers †", "55", "24", "4", "4", "229.16")
    22 = @(" Dubai (DSC)", " October 17 2020", "RCB won by 7 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Rajasthan Royals", "AB de Villiers †", "55", "22", "1", "6", "250.00")
    23 = @(" Sharjah", " October 12 2020", "RCB won by 82 runs", "Royal Challengers Bangalore", "Kolkata Knight Riders", "AB de Villiers †", "73", "33", "5", "6", "221.21")
    24 = @(" Dubai (DSC)", " October 10 2020", "RCB won by 37 runs", "Royal Challengers Bangalore", "Chennai Super Kings", "AB de Villiers †", "0", "2", "0", "0", "0.00")
    25 = @(" Abu Dhabi", " October 03 2020", "RCB won by 8 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Rajasthan Royals", "AB de Villiers †", "12", "10", "1", "0", "120.00")
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K")

foreach ($r in ($newRows.Keys | Sort-Object)) {
    $vals = $newRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + $r)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$i]
    }
}

Write-Output "done"
